$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns are treated as plain text,
# matching the source data which stores values such as "1.010" or
# "0.07619" as literal strings (so Excel does not coerce them to numbers
# and strip significant trailing zeros).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.072.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.823.15'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.54'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.009'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4656'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.71%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3639'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07299'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8675'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07619'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.28%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.858.68'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.83'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.473'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.011'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008641'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.474.46'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.198'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.57'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.089.86'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.80'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.865'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.27'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.103'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.52%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.19'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.092'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08930'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.958'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7319'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.455'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.141'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.35%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.526'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.20%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.074'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.95%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05270'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.88%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01917'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.45%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.938'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.123'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5216'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1633'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.258'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4858'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.010'
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '103.77'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.42%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.11'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.643'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06259'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.11%  '
